# Refresh cryptos snapshot: prices (col D) and 1h volume-% (col E) for most
# rows, plus rows 12/13 (Toncoin <-> Cardano) swapping their Coin/Link/Price/
# Volume data -- matching commit 'Updated cryptos list on Sun Oct 20 22:08:19
# UTC 2024 with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: a few Price cells (col D) are plain-looking decimals (e.g. '164.90',
# '1.00'). Assigning those bare would make Excel coerce them to a Number and
# silently lose the trailing zero (164.90 -> 164.9). Prefixing with a literal
# leading apostrophe forces Excel to keep them as text, same as typing
# '164.90 into the cell.

$ws.Range('D2').Value = '68.979.30'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '2.733.75'
$ws.Range('E3').Value = '  +3.37%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'602.21"
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').Value = "'164.90"
$ws.Range('E6').Value = '  +4.84%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'0.545"
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range('D9').Value = '2.731.50'
$ws.Range('E9').Value = '  +3.33%  '
$ws.Range('D10').Value = "'0.142"
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').Value = "'0.364"
$ws.Range('E12').Value = '  +3.62%  '
$ws.Range('B13').Value = 'Toncoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D13').Value = "'5.33"
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('D14').Value = "'28.69"
$ws.Range('E14').Value = '  +2.32%  '
$ws.Range('D15').Value = '3.234.20'
$ws.Range('E15').Value = '  +3.55%  '
$ws.Range('D16').Value = "'0.0000190"
$ws.Range('E16').Value = '  +0.90%  '
$ws.Range('D17').Value = '68.815.27'
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').Value = '2.727.76'
$ws.Range('E18').Value = '  +3.08%  '
$ws.Range('D19').Value = "'11.94"
$ws.Range('E19').Value = '  +4.96%  '
$ws.Range('D20').Value = "'7.73"
$ws.Range('E20').Value = '  +5.63%  '
$ws.Range('D21').Value = "'367.12"
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').Value = "'4.57"
$ws.Range('E22').Value = '  +3.52%  '
$ws.Range('D23').Value = "'4.95"
$ws.Range('E23').Value = '  +3.32%  '
$ws.Range('D24').Value = "'2.13"
$ws.Range('E24').Value = '  +3.52%  '
$ws.Range('D25').Value = "'73.98"
$ws.Range('E25').Value = '  -1.61%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = "'9.97"
$ws.Range('E27').Value = '  +2.67%  '
$ws.Range('D28').Value = '2.853.36'
$ws.Range('E28').Value = '  +2.72%  '
$ws.Range('D29').Value = "'0.0000106"
$ws.Range('E29').Value = '  +2.58%  '
$ws.Range('D30').Value = "'596.99"
$ws.Range('E30').Value = '  +7.19%  '
$ws.Range('D31').Value = "'0.996"
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('D32').Value = "'8.31"
$ws.Range('E32').Value = '  +3.76%  '
$ws.Range('D33').Value = "'1.46"
$ws.Range('E33').Value = '  +4.29%  '
$ws.Range('E34').Value = '  +5.74%  '
$ws.Range('D35').Value = "'0.133"
$ws.Range('E35').Value = '  +3.31%  '
$ws.Range('D36').Value = "'1.64"
$ws.Range('E36').Value = '  +5.58%  '
$ws.Range('D37').Value = "'0.999"
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = "'162.50"
$ws.Range('E38').Value = '  +1.51%  '
$ws.Range('D39').Value = "'20.06"
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('D40').Value = "'0.383"
$ws.Range('E40').Value = '  +3.22%  '
$ws.Range('D41').Value = "'1.93"
$ws.Range('E41').Value = '  +2.96%  '
$ws.Range('D42').Value = "'5.49"
$ws.Range('E42').Value = '  +3.10%  '
$ws.Range('D43').Value = "'2.72"
$ws.Range('E43').Value = '  +4.72%  '
$ws.Range('D44').Value = "'18.02"
$ws.Range('E44').Value = '  +1.27%  '
$ws.Range('D46').Value = '0.0₆0317'
$ws.Range('E46').Value = '  -4.24%  '
$ws.Range('D47').Value = "'158.77"
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').Value = "'3.96"
$ws.Range('E48').Value = '  +5.89%  '
$ws.Range('D49').Value = "'1.80"
$ws.Range('E49').Value = '  +6.75%  '
$ws.Range('D50').Value = "'0.612"
$ws.Range('E50').Value = '  +8.28%  '
$ws.Range('D51').Value = "'22.20"
$ws.Range('E51').Value = '  +0.57%  '
